$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.303865551948547
$ws.Range("B1").Value = 15
$ws.Range("D1").Value = 1.360026001930237
$ws.Range("E1").Value = 0.8257433772087097
